$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7454588356058228
$ws.Range("D2").Value = 0.4638850819412634

$ws.Range("C3").Value = 0.1379965594943917
$ws.Range("D3").Value = 0.8914984107874462

$ws.Range("C4").Value = -0.6654012230363781
$ws.Range("D4").Value = 0.5127102621285324

$ws.Range("C5").Value = 2.154753016767805
$ws.Range("D5").Value = 0.0423913271552272
$ws.Range("G5").Value = "Sí"

$ws.Range("C6").Value = -0.4910659922433831
$ws.Range("D6").Value = 0.6282427563702213

$ws.Range("C7").Value = -1.521839663027371
$ws.Range("D7").Value = 0.1422939299156285

$ws.Range("C8").Value = 1.397987542814258
$ws.Range("D8").Value = 0.1760538912425491

$ws.Range("C9").Value = -0.9248376023725877
$ws.Range("D9").Value = 0.3650843603587832

$ws.Range("C10").Value = 1.562534621962491
$ws.Range("D10").Value = 0.1324333835041722
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = 2.02313517055498
$ws.Range("D11").Value = 0.05537712452840715
$ws.Range("G11").Value = "No"
